# Re-generated CAM/BOM/CPL data: KiCad re-exported the top-side placement
# CSV and the workbook's query tables were refreshed, updating the
# position/rotation for R6 (row 29) and U1 (row 30).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B29").Value = 143.256
$ws.Range("C29").Value = -103.124
$ws.Range("E29").Value = 0

$ws.Range("B30").Value = 143.256
$ws.Range("C30").Value = -101.854
$ws.Range("E30").Value = 0

$ws.Range("E31").Select()
